# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late"/"Waived"/"Outstanding") columns to
#   O/P/Q.
# - Give the freshly inserted column the same display width as its
#   neighbour.
# - Switch the active sheet/tab from "NewLoanInput" to "Repayment schedule"
#   and leave the selection on R8 (one cell past the new last column), as
#   last recorded by the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N (pushes Late/Waived/Outstanding
# one column to the right).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.140625

# Make "Repayment schedule" the active sheet/tab and park the selection on
# R8, matching the recorded workbook state.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
